$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7173645496368408
$ws.Range("B1").Value = 1.028411746025085
$ws.Range("C1").Value = 1.959414482116699
$ws.Range("D1").Value = 3.384492874145508
$ws.Range("E1").Value = 3.6192946434021
